# Scheduled runner: refresh cached Universalis price snapshots and the
# derived Leve profit columns (H..N) across all job sheets.

$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param(
        [string]$SheetName,
        [string]$CellRef,
        $Value
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).Value = $Value
}

function Clear-CellValue {
    param(
        [string]$SheetName,
        [string]$CellRef
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).ClearContents()
}

# ---------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------
Set-CellValue "ALC" "H70" 20959282
Set-CellValue "ALC" "I70" 23953250
Set-CellValue "ALC" "J70" 1500
Set-CellValue "ALC" "K70" 71859750
Set-CellValue "ALC" "L70" 4500
Set-CellValue "ALC" "M70" -71859480
Set-CellValue "ALC" "N70" -5040

Set-CellValue "ALC" "H73" 20959282
Set-CellValue "ALC" "I73" 23953250
Set-CellValue "ALC" "J73" 1500
Set-CellValue "ALC" "K73" 71859750
Set-CellValue "ALC" "L73" 4500
Set-CellValue "ALC" "M73" -71858814
Set-CellValue "ALC" "N73" -6372

Set-CellValue "ALC" "H125" 3855.4707
Set-CellValue "ALC" "J125" 3516.2666
Set-CellValue "ALC" "L125" 31646.3994
Set-CellValue "ALC" "N125" -36566.39939999999

Set-CellValue "ALC" "H129" 1150.08
Set-CellValue "ALC" "I129" 441.66666
Set-CellValue "ALC" "K129" 1324.99998
Set-CellValue "ALC" "M129" 3675.00002

Set-CellValue "ALC" "H131" 7388.773
Set-CellValue "ALC" "J131" 10074
Set-CellValue "ALC" "L131" 30222
Set-CellValue "ALC" "N131" -40302

Set-CellValue "ALC" "H137" 1484.2667
Set-CellValue "ALC" "I137" 1484.2667
Set-CellValue "ALC" "J137" 0
Set-CellValue "ALC" "K137" 4452.800099999999
Set-CellValue "ALC" "L137" 0
Set-CellValue "ALC" "M137" -1902.800099999999
Clear-CellValue "ALC" "N137"

# ---------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------
Set-CellValue "ARM" "H32" 12015.477
Set-CellValue "ARM" "I32" 12946.732
Set-CellValue "ARM" "K32" 12946.732
Set-CellValue "ARM" "M32" -12659.732

Set-CellValue "ARM" "H92" 134820
Set-CellValue "ARM" "J92" 134820
Set-CellValue "ARM" "L92" 134820
Set-CellValue "ARM" "N92" -139812

# ---------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------
Set-CellValue "BSM" "H64" 2086
Set-CellValue "BSM" "I64" 2198.3333
Set-CellValue "BSM" "J64" 2001.75
Set-CellValue "BSM" "K64" 2198.3333
Set-CellValue "BSM" "L64" 2001.75
Set-CellValue "BSM" "M64" -1973.3333
Set-CellValue "BSM" "N64" -2451.75

Set-CellValue "BSM" "H67" 2086
Set-CellValue "BSM" "I67" 2198.3333
Set-CellValue "BSM" "J67" 2001.75
Set-CellValue "BSM" "K67" 2198.3333
Set-CellValue "BSM" "L67" 2001.75
Set-CellValue "BSM" "M67" -1418.3333
Set-CellValue "BSM" "N67" -3561.75

Set-CellValue "BSM" "H92" 262695.5
Set-CellValue "BSM" "J92" 262695.5
Set-CellValue "BSM" "L92" 262695.5
Set-CellValue "BSM" "N92" -267687.5

Set-CellValue "BSM" "H95" 0
Set-CellValue "BSM" "J95" 0
Set-CellValue "BSM" "L95" 0
Clear-CellValue "BSM" "N95"

# ---------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------
Set-CellValue "CRP" "H31" 3559.9473
Set-CellValue "CRP" "I31" 2123.923
Set-CellValue "CRP" "J31" 6671.3335
Set-CellValue "CRP" "K31" 2123.923
Set-CellValue "CRP" "L31" 6671.3335
Set-CellValue "CRP" "M31" -1828.923
Set-CellValue "CRP" "N31" -7261.3335

Set-CellValue "CRP" "H34" 3559.9473
Set-CellValue "CRP" "I34" 2123.923
Set-CellValue "CRP" "J34" 6671.3335
Set-CellValue "CRP" "K34" 2123.923
Set-CellValue "CRP" "L34" 6671.3335
Set-CellValue "CRP" "M34" -1921.923
Set-CellValue "CRP" "N34" -7075.3335

# ---------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------
Set-CellValue "CUL" "H106" 7325.5557
Set-CellValue "CUL" "J106" 7325.5557
Set-CellValue "CUL" "L106" 21976.6671
Set-CellValue "CUL" "N106" -23868.6671

Set-CellValue "CUL" "H109" 16644.857
Set-CellValue "CUL" "I109" 50427
Set-CellValue "CUL" "J109" 3132
Set-CellValue "CUL" "K109" 151281
Set-CellValue "CUL" "L109" 9396
Set-CellValue "CUL" "M109" -150241
Set-CellValue "CUL" "N109" -11476

Set-CellValue "CUL" "H112" 4988.778
Set-CellValue "CUL" "I112" 3499.5
Set-CellValue "CUL" "J112" 5414.2856
Set-CellValue "CUL" "K112" 10498.5
Set-CellValue "CUL" "L112" 16242.8568
Set-CellValue "CUL" "M112" -9390.5
Set-CellValue "CUL" "N112" -18458.8568

Set-CellValue "CUL" "H131" 16394727
Set-CellValue "CUL" "I131" 1730
Set-CellValue "CUL" "J131" 16950422
Set-CellValue "CUL" "K131" 5190
Set-CellValue "CUL" "L131" 50851266
Set-CellValue "CUL" "M131" -150
Set-CellValue "CUL" "N131" -50861346

Set-CellValue "CUL" "H133" 5292.231
Set-CellValue "CUL" "I133" 2576
Set-CellValue "CUL" "J133" 6989.875
Set-CellValue "CUL" "K133" 7728
Set-CellValue "CUL" "L133" 20969.625
Set-CellValue "CUL" "M133" -2668
Set-CellValue "CUL" "N133" -31089.625

# ---------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------
Set-CellValue "GSM" "H113" 2382.7585
Set-CellValue "GSM" "I113" 1317.5834
Set-CellValue "GSM" "J113" 3134.647
Set-CellValue "GSM" "K113" 1317.5834
Set-CellValue "GSM" "L113" 3134.647
Set-CellValue "GSM" "M113" 852.4166
Set-CellValue "GSM" "N113" -7474.647

Set-CellValue "GSM" "H123" 8972.529
Set-CellValue "GSM" "J123" 8972.529
Set-CellValue "GSM" "L123" 8972.529
Set-CellValue "GSM" "N123" -13872.529

Set-CellValue "GSM" "H126" 4026.6
Set-CellValue "GSM" "I126" 3953.8462
Set-CellValue "GSM" "J126" 4499.5
Set-CellValue "GSM" "K126" 11861.5386
Set-CellValue "GSM" "L126" 13498.5
Set-CellValue "GSM" "M126" -9391.5386
Set-CellValue "GSM" "N126" -18438.5

# ---------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------
Set-CellValue "LTW" "H68" 2666.6667
Set-CellValue "LTW" "I68" 2615.3845
Set-CellValue "LTW" "J68" 3000
Set-CellValue "LTW" "K68" 2615.3845
Set-CellValue "LTW" "L68" 3000
Set-CellValue "LTW" "M68" -1866.3845
Set-CellValue "LTW" "N68" -4498

Set-CellValue "LTW" "H71" 2666.6667
Set-CellValue "LTW" "I71" 2615.3845
Set-CellValue "LTW" "J71" 3000
Set-CellValue "LTW" "K71" 13076.9225
Set-CellValue "LTW" "L71" 15000
Set-CellValue "LTW" "M71" -9332.922500000001
Set-CellValue "LTW" "N71" -22488

Set-CellValue "LTW" "H94" 0
Set-CellValue "LTW" "J94" 0
Set-CellValue "LTW" "L94" 0
Clear-CellValue "LTW" "N94"

Set-CellValue "LTW" "H112" 30000
Set-CellValue "LTW" "J112" 30000
Set-CellValue "LTW" "L112" 30000
Set-CellValue "LTW" "N112" -32954

# ---------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------
Set-CellValue "WVR" "H2" 1000
Set-CellValue "WVR" "J2" 1000
Set-CellValue "WVR" "L2" 1000
Set-CellValue "WVR" "N2" -1224

Set-CellValue "WVR" "H21" 50000
Set-CellValue "WVR" "J21" 0
Set-CellValue "WVR" "L21" 0
Clear-CellValue "WVR" "N21"

Set-CellValue "WVR" "H35" 50000
Set-CellValue "WVR" "J35" 0
Set-CellValue "WVR" "L35" 0
Clear-CellValue "WVR" "N35"

Set-CellValue "WVR" "H101" 7278.4
Set-CellValue "WVR" "J101" 7278.4
Set-CellValue "WVR" "L101" 7278.4
Set-CellValue "WVR" "N101" -13768.4

Set-CellValue "WVR" "H113" 374.15
Set-CellValue "WVR" "I113" 476.7143
Set-CellValue "WVR" "J113" 318.92307
Set-CellValue "WVR" "K113" 1430.1429
Set-CellValue "WVR" "L113" 956.7692099999999
Set-CellValue "WVR" "M113" 739.8571000000002
Set-CellValue "WVR" "N113" -5296.76921
